{"js": "// Update the division problems in the table to the new set of values.\n// Each old \"A\u00f7B=\" text is replaced with its corresponding new \"A\u00f7B=\" text.\nconst replacements = [\n  [\"524\u00f79=\", \"266\u00f76=\"],\n  [\"612\u00f78=\", \"402\u00f77=\"],\n  [\"339\u00f72=\", \"822\u00f77=\"],\n  [\"972\u00f76=\", \"719\u00f77=\"],\n  [\"444\u00f78=\", \"507\u00f76=\"],\n  [\"878\u00f75=\", \"676\u00f78=\"],\n  [\"364\u00f78=\", \"335\u00f72=\"],\n  [\"662\u00f72=\", \"958\u00f77=\"],\n  [\"529\u00f78=\", \"133\u00f75=\"],\n  [\"227\u00f77=\", \"278\u00f78=\"],\n  [\"513\u00f79=\", \"920\u00f74=\"],\n  [\"657\u00f73=\", \"637\u00f75=\"],\n  [\"148\u00f75=\", \"293\u00f74=\"],\n  [\"356\u00f75=\", \"616\u00f77=\"],\n  [\"836\u00f77=\", \"193\u00f77=\"],\n  [\"364\u00f76=\", \"384\u00f77=\"],\n  [\"930\u00f78=\", \"804\u00f72=\"],\n  [\"521\u00f75=\", \"125\u00f76=\"],\n  [\"519\u00f79=\", \"917\u00f75=\"],\n  [\"589\u00f74=\", \"757\u00f75=\"],\n  [\"274\u00f73=\", \"976\u00f79=\"],\n  [\"499\u00f73=\", \"851\u00f77=\"],\n  [\"689\u00f79=\", \"608\u00f73=\"],\n  [\"236\u00f72=\", \"202\u00f78=\"],\n  [\"742\u00f72=\", \"920\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the table to the new set of values.\n# Each old \"A\u00f7B=\" text is replaced with its corresponding new \"A\u00f7B=\" text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"524\u00f79=\"; New = \"266\u00f76=\" },\n    @{ Old = \"612\u00f78=\"; New = \"402\u00f77=\" },\n    @{ Old = \"339\u00f72=\"; New = \"822\u00f77=\" },\n    @{ Old = \"972\u00f76=\"; New = \"719\u00f77=\" },\n    @{ Old = \"444\u00f78=\"; New = \"507\u00f76=\" },\n    @{ Old = \"878\u00f75=\"; New = \"676\u00f78=\" },\n    @{ Old = \"364\u00f78=\"; New = \"335\u00f72=\" },\n    @{ Old = \"662\u00f72=\"; New = \"958\u00f77=\" },\n    @{ Old = \"529\u00f78=\"; New = \"133\u00f75=\" },\n    @{ Old = \"227\u00f77=\"; New = \"278\u00f78=\" },\n    @{ Old = \"513\u00f79=\"; New = \"920\u00f74=\" },\n    @{ Old = \"657\u00f73=\"; New = \"637\u00f75=\" },\n    @{ Old = \"148\u00f75=\"; New = \"293\u00f74=\" },\n    @{ Old = \"356\u00f75=\"; New = \"616\u00f77=\" },\n    @{ Old = \"836\u00f77=\"; New = \"193\u00f77=\" },\n    @{ Old = \"364\u00f76=\"; New = \"384\u00f77=\" },\n    @{ Old = \"930\u00f78=\"; New = \"804\u00f72=\" },\n    @{ Old = \"521\u00f75=\"; New = \"125\u00f76=\" },\n    @{ Old = \"519\u00f79=\"; New = \"917\u00f75=\" },\n    @{ Old = \"589\u00f74=\"; New = \"757\u00f75=\" },\n    @{ Old = \"274\u00f73=\"; New = \"976\u00f79=\" },\n    @{ Old = \"499\u00f73=\"; New = \"851\u00f77=\" },\n    @{ Old = \"689\u00f79=\"; New = \"608\u00f73=\" },\n    @{ Old = \"236\u00f72=\"; New = \"202\u00f78=\" },\n    @{ Old = \"742\u00f72=\"; New = \"920\u00f77=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.Old,   # FindText\n        $true,    # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # Wrap (wdFindContinue)\n        $false,   # Format\n        $r.New,   # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n$d.Save()\n"}
